$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: store a literal text value (e.g. "24.1%") in a cell without Excel's
# automatic "looks like a percentage" type coercion turning it into a number,
# and without provisioning a brand-new cell style (NumberFormat/quote-prefix
# tricks both end up creating a new xf record). Routing it through a formula
# and then pasting-values collapses it back down to a plain shared-string
# cell that keeps the range's pre-existing style.
function Set-LiteralText($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# G2 - "Recorded By" list reordered
$ws.Cells.Item(2, 7).Value = "System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# G3 - "Recorded By" list reordered
$ws.Cells.Item(3, 7).Value = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# L6 - Recorded Sessions: 6 -> 7
$ws.Cells.Item(6, 12).Value = 7

# L7 - Missing Sessions: 1 -> 0
$ws.Cells.Item(7, 12).Value = 0

# G9 - "Recorded By" list reordered
$ws.Cells.Item(9, 7).Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# L9 - Coverage %: 20.7% -> 24.1%
Set-LiteralText $ws.Cells.Item(9, 12) "24.1%"

# L10 - Average Attendance %: 28.0% -> 27.1%
Set-LiteralText $ws.Cells.Item(10, 12) "27.1%"

# Column I width: 14 -> 10 (match column H, which is already width 10 -
# copying the property value keeps the exact internal/XML representation
# instead of hand-rounding a "10" that Excel's char-width<->XML-width
# conversion would otherwise turn into 10.8333...)
$ws.Columns(9).ColumnWidth = $ws.Columns(8).ColumnWidth

# Row 15 - PARASITOLOGY / C1 / session 2 has now been recorded (was Pending)
# Match the formatting of the other "Recorded" rows (e.g. row 9) instead of
# the "Not Recorded" style it had before.
$src = $ws.Range("A9:I9")
$dst = $ws.Range("A15:I15")
$src.Copy()
$dst.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 7).Value = "Rania.a.youssef@med.asu.edu.eg"
$ws.Cells.Item(15, 8).Value = "56/251"
$ws.Cells.Item(15, 9).Value = "Recorded"

# Row 15 rolled-up statistics columns (O, P, R, S)
$ws.Cells.Item(15, 15).Value = 7
$ws.Cells.Item(15, 16).Value = 0
Set-LiteralText $ws.Cells.Item(15, 18) "24.1%"
Set-LiteralText $ws.Cells.Item(15, 19) "27.1%"

# G28 - "Recorded By" list reordered
$ws.Cells.Item(28, 7).Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
